# Updates cryptos list price/volume data (and two rank-swap rows) per the
# Fri Aug 25 14:23:42 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Column D holds price strings (e.g. "26.178.10") that Excel would
    # otherwise auto-coerce into a number; force text, write, then drop
    # back to the Normal style so no stray formatting is left behind.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '26.178.10'
$ws.Range('E2').Value = '  -0.79%  '
# Row 3
Set-TextValue 'D3' '1.664.28'
$ws.Range('E3').Value = '  -0.12%  '
# Row 4
$ws.Range('E4').Value = '  -0.18%  '
# Row 5
Set-TextValue 'D5' '217.83'
$ws.Range('E5').Value = '  -0.99%  '
# Row 6
Set-TextValue 'D6' '0.5260'
$ws.Range('E6').Value = '  +0.20%  '
# Row 7
Set-TextValue 'D7' '1.002'
$ws.Range('E7').Value = '  -0.19%  '
# Row 8
Set-TextValue 'D8' '0.2645'
$ws.Range('E8').Value = '  -0.96%  '
# Row 9
Set-TextValue 'D9' '0.06286'
$ws.Range('E9').Value = '  -1.10%  '
# Row 10
Set-TextValue 'D10' '20.78'
$ws.Range('E10').Value = '  -3.91%  '
# Row 11
Set-TextValue 'D11' '0.07755'
$ws.Range('E11').Value = '  -0.12%  '
# Row 12
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D12' '4.467'
$ws.Range('E12').Value = '  +0.09%  '
# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D13' '1.633.15'
$ws.Range('E13').Value = '  -2.00%  '
# Row 14
Set-TextValue 'D14' '1.887.54'
$ws.Range('E14').Value = '  -0.36%  '
# Row 15
Set-TextValue 'D15' '0.5478'
$ws.Range('E15').Value = '  -0.67%  '
# Row 16
Set-TextValue 'D16' '0.0₅8123'
$ws.Range('E16').Value = '  -1.72%  '
# Row 17
Set-TextValue 'D17' '64.94'
$ws.Range('E17').Value = '  -0.78%  '
# Row 18
Set-TextValue 'D18' '26.180.00'
$ws.Range('E18').Value = '  -0.85%  '
# Row 19
$ws.Range('E19').Value = '  -0.11%  '
# Row 20
Set-TextValue 'D20' '4.594'
$ws.Range('E20').Value = '  -2.84%  '
# Row 21
Set-TextValue 'D21' '192.15'
$ws.Range('E21').Value = '  -0.56%  '
# Row 22
Set-TextValue 'D22' '10.04'
$ws.Range('E22').Value = '  -2.14%  '
# Row 23
Set-TextValue 'D23' '6.014'
$ws.Range('E23').Value = '  -3.88%  '
# Row 24
Set-TextValue 'D24' '1.003'
$ws.Range('E24').Value = '  -0.26%  '
# Row 25
Set-TextValue 'D25' '137.57'
$ws.Range('E25').Value = '  -1.00%  '
# Row 26
Set-TextValue 'D26' '0.1240'
$ws.Range('E26').Value = '  -1.88%  '
# Row 27
Set-TextValue 'D27' '7.258'
$ws.Range('E27').Value = '  -1.55%  '
# Row 28
Set-TextValue 'D28' '16.22'
$ws.Range('E28').Value = '  +0.25%  '
# Row 29
Set-TextValue 'D29' '1.402'
$ws.Range('E29').Value = '  -1.04%  '
# Row 30
Set-TextValue 'D30' '0.05981'
$ws.Range('E30').Value = '  -2.00%  '
# Row 31
Set-TextValue 'D31' '1.279'
$ws.Range('E31').Value = '  -0.89%  '
# Row 32
Set-TextValue 'D32' '3.531'
$ws.Range('E32').Value = '  -1.39%  '
# Row 33
Set-TextValue 'D33' '3.268'
$ws.Range('E33').Value = '  -3.68%  '
# Row 34
Set-TextValue 'D34' '1.581'
$ws.Range('E34').Value = '  -5.63%  '
# Row 35
Set-TextValue 'D35' '0.9615'
$ws.Range('E35').Value = '  -3.68%  '
# Row 36
Set-TextValue 'D36' '2.417'
$ws.Range('E36').Value = '  -0.24%  '
# Row 37
$ws.Range('E37').Value = '  -0.38%  '
# Row 38
Set-TextValue 'D38' '0.5669'
# Row 39
Set-TextValue 'D39' '0.01599'
$ws.Range('E39').Value = '  -0.55%  '
# Row 40
Set-TextValue 'D40' '5.928'
$ws.Range('E40').Value = '  -1.10%  '
# Row 41
Set-TextValue 'D41' '0.8512'
$ws.Range('E41').Value = '  -0.51%  '
# Row 42
$ws.Range('E42').Value = '  -0.15%  '
# Row 43
Set-TextValue 'D43' '100.89'
$ws.Range('E43').Value = '  +0.39%  '
# Row 44
Set-TextValue 'D44' '1.008.31'
$ws.Range('E44').Value = '  -6.92%  '
# Row 45
Set-TextValue 'D45' '1.801.41'
$ws.Range('E45').Value = '  -0.48%  '
# Row 46
Set-TextValue 'D46' '56.85'
$ws.Range('E46').Value = '  -1.61%  '
# Row 47
$ws.Range('E47').Value = '  -3.75%  '
# Row 48
Set-TextValue 'D48' '1.002'
$ws.Range('E48').Value = '  +0.08%  '
# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D49' '8.015'
$ws.Range('E49').Value = '  -1.23%  '
# Row 50
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D50' '0.4326'
$ws.Range('E50').Value = '  +2.21%  '
# Row 51
Set-TextValue 'D51' '0.05148'
$ws.Range('E51').Value = '  -1.09%  '
